# Remove the unused "Global" sheet and rename the header of the
# "服务类别" (service category) sheet from "服务类别名称" to "值".

$wb = $excel.ActiveWorkbook

# Rename the header cell before removing the other sheet.
$ws = $wb.Worksheets.Item("服务类别")
$ws.Range("A1").Value = "值"

# Delete the "Global" sheet entirely (suppress the confirmation alert).
$excel.DisplayAlerts = $false
$global = $wb.Worksheets.Item("Global")
$global.Delete()
